$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.895.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.124.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.17%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '617.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.11'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.393'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.122.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.755'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.205'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.39%  '

$ws.Range("E15").Value = '  +1.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.665.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.42%  '

$ws.Range("E17").Value = '  +0.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.127.04'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '456.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000204'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '89.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +46.56%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.76%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  +19.68%  '

$ws.Range("E31").Value = '  -0.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.229'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.41%  '

$ws.Range("E33").Value = '  -7.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.37'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.16%  '

$ws.Range("E35").Value = '  +5.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.91%  '

$ws.Range("B37").Value = 'PancakeSwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.12%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.52'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '493.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.35%  '

$ws.Range("E40").Value = '  +0.67%  '

$ws.Range("E41").Value = '  -7.65%  '

$ws.Range("E42").Value = '  +3.53%  '

$ws.Range("E43").Value = '  -6.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.20'
$ws.Range("D44").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.708'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.01%  '

$ws.Range("E47").Value = '  -1.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '156.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.41%  '

$ws.Range("E49").Value = '  -0.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.49'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.07%  '

$ws.Range("E51").Value = '  +1.32%  '
